$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Chicken "
$ws.Range("B2").Value = 32532.0
$ws.Range("C2").Value = "C"
$ws.Range("D2").Value = 1.0
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 220.56
$ws.Range("G2").Value = 60.51
$ws.Range("H2").Value = 1000.0
$ws.Range("I2").Value = 100.0
